# Fill in the second progress table (rows 12-15) with the team's
# completed-task information, and update the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A12").Value = "余舒章"
$ws.Range("B12").Value = "修改完善pc端用例"
$ws.Range("C12").Value = "已完成"

$ws.Range("A13").Value = "王嘉宇"
$ws.Range("B13").Value = "修改完善pc端用例"
$ws.Range("C13").Value = "已完成"

$ws.Range("A14").Value = "许俊杰"
$ws.Range("B14").Value = "修改完善android端用例"
$ws.Range("C14").Value = "已完成"

$ws.Range("A15").Value = "庞森杰"
$ws.Range("B15").Value = "修改完善android端用例"
$ws.Range("C15").Value = "已完成"

# Update the saved selection to match the author's final cursor position.
[void]$ws.Range("F11").Select()
